# "Fixed StudyComb for Faceted Filters ICDC"
# The StatQuery column (C2:C4) held a stale/broken Cypher aggregate query.
# Replace it with the corrected query text (same text repeated in all three
# data-tab rows: CasesTab / SamplesTab / FilesTab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStatQuery = "MATCH (demo:demographic)`nWHERE demo.breed IN ['Australian Shepherd']`nMATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)`nOPTIONAL MATCH (c)<-[*]-(samp:sample)`nOPTIONAL MATCH (c)<-[*]-(f:file)`nRETURN `n`tcount(DISTINCT(f)) as number_of_files, `n`tcount(DISTINCT(samp)) as number_of_sample, `n`tcount(DISTINCT(c)) as number_of_cases, `n`tcount(DISTINCT(s)) as number_of_study"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the view state left behind by the save (scrolled/zoomed to row 4,
# selection on B4, zoom reset to 100%).
$ws.Range("B4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
